$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Employee" label into C7 (next to "Applicant" in B7)
$ws.Range("C7").Value = "Employee"

# Update the active selection to D5 (approving submission moves focus there)
$ws.Range("D5").Select()

# Force recalculation so TODAY() formulas refresh their cached values
$excel.Calculate()
